$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Measured Data")

# Insert a new column before column F ("Bleeding Ratio, Beta"), shifting
# Flowrate/Viable/Dead/Total Cell Concentration columns one slot to the right.
$ws.Columns("F:F").Insert()

# New header label for the inserted column.
$ws.Range("F2").Value = "Bleeding Ratio, Beta"

# Bleeding ratio values for each data row (mostly 0, with two exceptions).
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0.2
$ws.Range("F20").Value = 0.05
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0

# Make "Measured Data" the active/selected sheet with F2 selected.
$ws.Select()
$ws.Range("F2").Select()
